$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header: "Datos actualizados" timestamp update
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 17:22"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 929028
$ws.Range("C4").Value = 3796
$ws.Range("D4").Value = 110504
$ws.Range("E4").Value = 766153
$ws.Range("G4").Value = 178
$ws.Range("H4").Value = 52371

# Row 17: Paises Bajos - Casos criticos updated
$ws.Range("F17").Value = 959

# Row 25: Arabia Saudita - Casos criticos updated
$ws.Range("F25").Value = 115

# Rows 29-31: Chile moves above Japon/Singapur in ranking (sorted by casos totales)
# Row 29 becomes Chile with refreshed data
$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 12858
$ws.Range("C29").Value = 552
$ws.Range("D29").Value = 6746
$ws.Range("E29").Value = 5931
$ws.Range("F29").Value = 418
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 181

# Row 30 becomes Japon (previous row 29 data, unchanged)
$ws.Range("A30").Value = "Japon"
$ws.Range("B30").Value = 12829
$ws.Range("C30").Value = 117
$ws.Range("D30").Value = 1530
$ws.Range("E30").Value = 10954
$ws.Range("F30").Value = 263
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 345

# Row 31 becomes Singapur (previous row 30 data, unchanged)
$ws.Range("A31").Value = "Singapur"
$ws.Range("B31").Value = 12693
$ws.Range("C31").Value = 618
$ws.Range("D31").Value = 956
$ws.Range("E31").Value = 11725
$ws.Range("F31").Value = 24
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 12

# Row 64: Grecia - refreshed data
$ws.Range("B64").Value = 2506
$ws.Range("C64").Value = 16
$ws.Range("E64").Value = 1799
$ws.Range("F64").Value = 47

# Rows 91-92: Republica de Chipre moves above Bolivia in ranking
# Row 91 becomes Republica de Chipre with refreshed data
$ws.Range("A91").Value = "Republica de Chipre"
$ws.Range("B91").Value = 810
$ws.Range("C91").Value = 6
$ws.Range("D91").Value = 98
$ws.Range("E91").Value = 698
$ws.Range("F91").Value = 15
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 14

# Row 92 becomes Bolivia (previous row 91 data, unchanged)
$ws.Range("A92").Value = "Bolivia"
$ws.Range("B92").Value = 807
$ws.Range("C92").Value = 104
$ws.Range("D92").Value = 54
$ws.Range("E92").Value = 709
$ws.Range("F92").Value = 3
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 44

# Row 118: Mauricio - refreshed data
$ws.Range("D118").Value = 295
$ws.Range("E118").Value = 27

# Row 144: Trinidad y Tobago - refreshed data
$ws.Range("D144").Value = 52
$ws.Range("E144").Value = 55
